# Applies the commit "Added barplots for Rdeep,bully (incomplete) and rand.
# R source files present" to the workbook:
#   1. Rename the shared string "ml_advancedSImple" -> "ml_enriched"
#      (every cell that used it keeps showing the corrected label).
#   2. Sheet "nieuwe namen": the Round-data block that was mistakenly left
#      on the "ml_stripped" row (row 32) really belongs to the "ml_minimal"
#      row (row 31) - move it there and clear row 32's numbers. (The
#      dependent AVERAGE()/... formulas in rows 41/42 recalc on their own.)
#   3. Sheet "nieuwe namen": add the missing "Against Rdeep" Mean / Part
#      won / Percentage won summary block (rows 73-79), mirroring the
#      existing "Against Rand" (rows 13-19) and "Against Bully"
#      (rows 40-46) blocks, but pointing at rows 64-69.
#   4. Sheet "Oude benamingen": same kind of mix-up as (2) but here the
#      numbers need to move from row 16 ("ml_minimal") to row 15
#      ("ml_stripped").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Shared-string fix: "ml_advancedSImple" -> "ml_enriched"
# ---------------------------------------------------------------------
$wsNew = $wb.Worksheets.Item("nieuwe namen")
$wsOld = $wb.Worksheets.Item("Oude benamingen")

foreach ($addr in @("A6", "B17", "A34", "B44", "A67")) {
    $wsNew.Range($addr).Value2 = "ml_enriched"
}

# ---------------------------------------------------------------------
# 2. "nieuwe namen": move the Round 1..9 numbers from row 32 to row 31
# ---------------------------------------------------------------------
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S")
foreach ($col in $cols) {
    $val = $wsNew.Range($col + "32").Value2
    $wsNew.Range($col + "31").Value2 = $val
}
$wsNew.Range("B32:S32").ClearContents()

# ---------------------------------------------------------------------
# 3. "nieuwe namen": add the "Against Rdeep" mean/part-won block
# ---------------------------------------------------------------------
$wsNew.Range("B73").Value2 = "Mean"
$wsNew.Range("C73").Value2 = "mean won games"
$wsNew.Range("E73").Value2 = "Part won"
$wsNew.Range("F73").Value2 = "Percentage won"

$labels = @{ 74 = "ml_minimal"; 75 = "ml_stripped"; 76 = "ml"; 77 = "ml_enriched"; 78 = "ml_advanced"; 79 = "ml_combined" }
$sourceRow = @{ 74 = 64; 75 = 65; 76 = 66; 77 = 67; 78 = 68; 79 = 69 }

foreach ($row in 74..79) {
    $src = $sourceRow[$row]
    $wsNew.Range("B$row").Value2 = $labels[$row]
    $wsNew.Range("C$row").Formula = "=AVERAGE(B$src,D$src,F$src,H$src,J$src,L$src,N$src,P$src,R$src)"
    $wsNew.Range("C$row").NumberFormat = "0.00"
    $wsNew.Range("D$row").NumberFormat = "0.00"
    $wsNew.Range("E$row").Formula = "=C$row/120"
    $wsNew.Range("E$row").NumberFormat = "0.00"
    $wsNew.Range("F$row").Formula = "=E$row*100"
    $wsNew.Range("F$row").NumberFormat = "0.00"
}

# ---------------------------------------------------------------------
# 4. "Oude benamingen": move the Round 1..9 numbers from row 16 to row 15
# ---------------------------------------------------------------------
foreach ($col in $cols) {
    $val = $wsOld.Range($col + "16").Value2
    $wsOld.Range($col + "15").Value2 = $val
}
$wsOld.Range("B16:S16").ClearContents()
